# Appending similar quantities (column names) to the appropriate subsection.
# Adds a new "sheet_3" worksheet at the end of the workbook holding the
# repeated column-name subsections (quantity_1 / quantity_2) together with
# their per-row values, and updates the view/selection state to match.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# --- Add the new worksheet after the last existing sheet ---------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "sheet_3"

# --- Populate sheet_3 with the repeated quantity subsections ------------
$headers = @("quantity_1", "quantity_1", "quantity_2", "quantity_2")
$row0 = @("q1_d0_r0", "q1_d1_r0", "q2_d0_r0", "q2_d1_r0")
$row1 = @("q1_d0_r1", "q1_d1_r1", "q2_d0_r1", "q2_d1_r1")

for ($c = 0; $c -lt 4; $c++) {
  $ws3.Cells.Item(1, $c + 1).Value = $headers[$c]
}
for ($c = 0; $c -lt 4; $c++) {
  $ws3.Cells.Item(2, $c + 1).Value = $row0[$c]
}
for ($c = 0; $c -lt 4; $c++) {
  $ws3.Cells.Item(3, $c + 1).Value = $row1[$c]
}

# Size the columns to fit their (short) content, like the authored sheet.
$ws3.Columns("A:D").AutoFit()

# --- View state: sheet_1 keeps a normal selection, sheet_3 becomes the
#     active/visible tab with its own selection -------------------------
$ws1.Activate() | Out-Null
$ws1.Range("C9").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("D7").Select() | Out-Null

# Scroll the workbook tabs so the first visible tab is sheet_2 (index 1),
# matching the authored workbook view.
$win = $wb.Windows.Item(1)
$win.ScrollWorkbookTabs(1, 1) | Out-Null
